$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update the header text for the "State ID" column to "State Alias"
$ws.Range("I1").Value = "State Alias"

# Update the selected cell to match the saved selection in the workbook
[void]$ws.Range("I1").Select()
